$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(
    "Rua Oswaldo de Oliveira Lima 575 – Parque Santa Rosa/Suzano",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Casa Branca/Suzano", 2
)

$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Execute(
    "CEP: 08663-310",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CEP: 08663-310 - SP", 2
)
